# Generate Report for Handback
#
# The "ee154e4d-69c6-475a-83f5-9caf1fb6f52e.md" file has been handed back
# and is now in sync with en-US, so the localization-status report needs
# to reflect the new status + timestamps, and the old "out of date" error
# needs to be cleared.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- "Overview" sheet: zh-cn / de-de status columns for the handed-back file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# --- "zh-cn" sheet: Status / Latest Handback DateTime / Error Detail ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("K3").Value = "2016-08-20 06:53:41"
$zhcn.Range("P3").Value = ""

# --- "de-de" sheet: Status / Latest Handback DateTime / Error Detail ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("K3").Value = "2016-08-20 06:53:47"
$dede.Range("P3").Value = ""
